$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header fields added to the shops form (pincode, village, tehsil)
# following the existing operatorName / address headers.
$ws.Range("C1").Value = "pincode"
$ws.Range("D1").Value = "village"
$ws.Range("E1").Value = "tehsil"

# Match the formatting already applied to the header row (font size 12,
# vertical-center alignment) so the new cells look like A1/B1/C1.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Move the active selection off the header row, onto F2.
$ws.Range("F2").Select() | Out-Null
